$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: extend with P1=14, Q1=15, matching the header style (s="1") ---
$ws.Range("O1").Copy($ws.Range("P1:Q1")) | Out-Null
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# --- Data rows 2-25: updated simulation results for columns B,C,D,E,G,H,I; O now 0; new P=0, Q=<result> ---
# Row 2
$ws.Cells.Item(2, 2).Value = 24.46292630714852
$ws.Cells.Item(2, 3).Value = 18.89699789604578
$ws.Cells.Item(2, 4).Value = 7.171401745091684
$ws.Cells.Item(2, 5).Value = 29.14883178299275
$ws.Cells.Item(2, 7).Value = 2.07333962946678
$ws.Cells.Item(2, 8).Value = 3.117916045888993
$ws.Cells.Item(2, 9).Value = 3.078024169146623
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 18.24709438119277

# Row 3
$ws.Cells.Item(3, 2).Value = 22.83994395951069
$ws.Cells.Item(3, 3).Value = 17.74351541351962
$ws.Cells.Item(3, 4).Value = 6.67375983682823
$ws.Cells.Item(3, 5).Value = 27.21992337864157
$ws.Cells.Item(3, 7).Value = 2.079191564016929
$ws.Cells.Item(3, 8).Value = 2.836422433677704
$ws.Cells.Item(3, 9).Value = 2.849172163516276
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = 0
$ws.Cells.Item(3, 17).Value = 17.74519587534337

# Row 4
$ws.Cells.Item(4, 2).Value = 21.78386351698823
$ws.Cells.Item(4, 3).Value = 17.00306022775515
$ws.Cells.Item(4, 4).Value = 6.351269465980063
$ws.Cells.Item(4, 5).Value = 25.97006093646953
$ws.Cells.Item(4, 7).Value = 2.082885124810285
$ws.Cells.Item(4, 8).Value = 2.658068073862267
$ws.Cells.Item(4, 9).Value = 2.705437727640331
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 17.43881980757267

# Row 5
$ws.Cells.Item(5, 2).Value = 21.33823227521352
$ws.Cells.Item(5, 3).Value = 16.70920078914408
$ws.Cells.Item(5, 4).Value = 6.219406939647969
$ws.Cells.Item(5, 5).Value = 25.44386509554195
$ws.Cells.Item(5, 7).Value = 2.084424831484577
$ws.Cells.Item(5, 8).Value = 2.583766853623442
$ws.Cells.Item(5, 9).Value = 2.646543087774031
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(5, 17).Value = 17.29927869223101

# Row 6
$ws.Cells.Item(6, 2).Value = 21.26306108347954
$ws.Cells.Item(6, 3).Value = 16.67941894030965
$ws.Cells.Item(6, 4).Value = 6.200016273343929
$ws.Cells.Item(6, 5).Value = 25.35510664747364
$ws.Cells.Item(6, 7).Value = 2.084692456847843
$ws.Cells.Item(6, 8).Value = 2.571100075244397
$ws.Cells.Item(6, 9).Value = 2.637546862487593
$ws.Cells.Item(6, 15).Value = 0
$ws.Cells.Item(6, 16).Value = 0
$ws.Cells.Item(6, 17).Value = 17.25752537730386

# Row 7
$ws.Cells.Item(7, 2).Value = 21.77720934598569
$ws.Cells.Item(7, 3).Value = 17.0512390098225
$ws.Cells.Item(7, 4).Value = 6.349368694163031
$ws.Cells.Item(7, 5).Value = 25.96200724157499
$ws.Cells.Item(7, 7).Value = 2.082933600040907
$ws.Cells.Item(7, 8).Value = 2.65642074116166
$ws.Cells.Item(7, 9).Value = 2.706495202457873
$ws.Cells.Item(7, 15).Value = 0
$ws.Cells.Item(7, 16).Value = 0
$ws.Cells.Item(7, 17).Value = 17.38654323465399

# Row 8
$ws.Cells.Item(8, 2).Value = 23.91475794785805
$ws.Cells.Item(8, 3).Value = 18.57013631689834
$ws.Cells.Item(8, 4).Value = 7.003185836426347
$ws.Cells.Item(8, 5).Value = 28.49601075291111
$ws.Cells.Item(8, 7).Value = 2.075371690766244
$ws.Cells.Item(8, 8).Value = 3.02116568231494
$ws.Cells.Item(8, 9).Value = 3.001456095152273
$ws.Cells.Item(8, 15).Value = 0
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(8, 17).Value = 18.00926780137992

# Row 9
$ws.Cells.Item(9, 2).Value = 27.63872934554203
$ws.Cells.Item(9, 3).Value = 21.20823741394084
$ws.Cells.Item(9, 4).Value = 8.152847446917766
$ws.Cells.Item(9, 5).Value = 32.95147110948562
$ws.Cells.Item(9, 7).Value = 2.06125623692083
$ws.Cells.Item(9, 8).Value = 3.695385720470106
$ws.Cells.Item(9, 9).Value = 3.554045888935708
$ws.Cells.Item(9, 15).Value = 0
$ws.Cells.Item(9, 16).Value = 0
$ws.Cells.Item(9, 17).Value = 19.32751093831113

# Row 10
$ws.Cells.Item(10, 2).Value = 30.08872485141643
$ws.Cells.Item(10, 3).Value = 23.00010264714325
$ws.Cells.Item(10, 4).Value = 8.917750377196473
$ws.Cells.Item(10, 5).Value = 35.03646330801931
$ws.Cells.Item(10, 7).Value = 2.051598387207084
$ws.Cells.Item(10, 8).Value = 4.128242123332592
$ws.Cells.Item(10, 9).Value = 3.93790457072418
$ws.Cells.Item(10, 15).Value = 0
$ws.Cells.Item(10, 16).Value = 0
$ws.Cells.Item(10, 17).Value = 20.06580119473863

# Row 11
$ws.Cells.Item(11, 2).Value = 31.10048171091787
$ws.Cells.Item(11, 3).Value = 23.75703352180395
$ws.Cells.Item(11, 4).Value = 9.237876338088819
$ws.Cells.Item(11, 5).Value = 28.4486128954523
$ws.Cells.Item(11, 7).Value = 2.049659633807174
$ws.Cells.Item(11, 8).Value = 4.496135585019219
$ws.Cells.Item(11, 9).Value = 4.032049582675228
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(11, 16).Value = 0
$ws.Cells.Item(11, 17).Value = 18.5199530720283

# Row 12
$ws.Cells.Item(12, 2).Value = 31.46113236288211
$ws.Cells.Item(12, 3).Value = 23.97530180258628
$ws.Cells.Item(12, 4).Value = 9.352745217663353
$ws.Cells.Item(12, 5).Value = 22.49800092346557
$ws.Cells.Item(12, 7).Value = 2.049724220782633
$ws.Cells.Item(12, 8).Value = 5.302163226731279
$ws.Cells.Item(12, 9).Value = 4.038467677042033
$ws.Cells.Item(12, 15).Value = 0
$ws.Cells.Item(12, 16).Value = 0
$ws.Cells.Item(12, 17).Value = 17.11277427944507

# Row 13
$ws.Cells.Item(13, 2).Value = 31.35466063464835
$ws.Cells.Item(13, 3).Value = 23.86149282393163
$ws.Cells.Item(13, 4).Value = 9.320188813591983
$ws.Cells.Item(13, 5).Value = 16.56310161362051
$ws.Cells.Item(13, 7).Value = 2.051365769174462
$ws.Cells.Item(13, 8).Value = 6.317728984374356
$ws.Cells.Item(13, 9).Value = 3.979873568238724
$ws.Cells.Item(13, 15).Value = 0
$ws.Cells.Item(13, 16).Value = 0
$ws.Cells.Item(13, 17).Value = 15.63412677740119

# Row 14
$ws.Cells.Item(14, 2).Value = 31.06815437475379
$ws.Cells.Item(14, 3).Value = 23.64372433656073
$ws.Cells.Item(14, 4).Value = 9.230612048821721
$ws.Cells.Item(14, 5).Value = 12.47452485418829
$ws.Cells.Item(14, 7).Value = 2.053198419658707
$ws.Cells.Item(14, 8).Value = 7.113124819939758
$ws.Cells.Item(14, 9).Value = 3.911245168130808
$ws.Cells.Item(14, 15).Value = 0
$ws.Cells.Item(14, 16).Value = 0
$ws.Cells.Item(14, 17).Value = 14.54976317536127

# Row 15
$ws.Cells.Item(15, 2).Value = 30.89671262648975
$ws.Cells.Item(15, 3).Value = 23.53432490604464
$ws.Cells.Item(15, 4).Value = 9.176874170981923
$ws.Cells.Item(15, 5).Value = 11.48676480888394
$ws.Cells.Item(15, 7).Value = 2.054025382061608
$ws.Cells.Item(15, 8).Value = 7.295430768928355
$ws.Cells.Item(15, 9).Value = 3.880343363681348
$ws.Cells.Item(15, 15).Value = 0
$ws.Cells.Item(15, 16).Value = 0
$ws.Cells.Item(15, 17).Value = 14.25574320816903

# Row 16
$ws.Cells.Item(16, 2).Value = 29.91491546860326
$ws.Cells.Item(16, 3).Value = 22.8455622040809
$ws.Cells.Item(16, 4).Value = 8.868949322270886
$ws.Cells.Item(16, 5).Value = 11.3119336879149
$ws.Cells.Item(16, 7).Value = 2.057784129828858
$ws.Cells.Item(16, 8).Value = 7.007537467995535
$ws.Cells.Item(16, 9).Value = 3.731126056603214
$ws.Cells.Item(16, 15).Value = 0
$ws.Cells.Item(16, 16).Value = 0
$ws.Cells.Item(16, 17).Value = 14.15859373800325

# Row 17
$ws.Cells.Item(17, 2).Value = 29.30381131141779
$ws.Cells.Item(17, 3).Value = 22.42833852179358
$ws.Cells.Item(17, 4).Value = 8.677498781340381
$ws.Cells.Item(17, 5).Value = 13.39918030100123
$ws.Cells.Item(17, 7).Value = 2.059726120998349
$ws.Cells.Item(17, 8).Value = 6.310778786467217
$ws.Cells.Item(17, 9).Value = 3.652297513878179
$ws.Cells.Item(17, 15).Value = 0
$ws.Cells.Item(17, 16).Value = 0
$ws.Cells.Item(17, 17).Value = 14.67091514957058

# Row 18
$ws.Cells.Item(18, 2).Value = 28.95853618611775
$ws.Cells.Item(18, 3).Value = 22.17307409356726
$ws.Cells.Item(18, 4).Value = 8.568780160373137
$ws.Cells.Item(18, 5).Value = 17.96953379334497
$ws.Cells.Item(18, 7).Value = 2.06016019333027
$ws.Cells.Item(18, 8).Value = 5.269627560738255
$ws.Cells.Item(18, 9).Value = 3.627323536352864
$ws.Cells.Item(18, 15).Value = 0
$ws.Cells.Item(18, 16).Value = 0
$ws.Cells.Item(18, 17).Value = 15.80152295421314

# Row 19
$ws.Cells.Item(19, 2).Value = 28.8600044156729
$ws.Cells.Item(19, 3).Value = 22.15708161533562
$ws.Cells.Item(19, 4).Value = 8.536600123578127
$ws.Cells.Item(19, 5).Value = 24.26567657739218
$ws.Cells.Item(19, 7).Value = 2.059173954332831
$ws.Cells.Item(19, 8).Value = 4.287378393464564
$ws.Cells.Item(19, 9).Value = 3.659909816288889
$ws.Cells.Item(19, 15).Value = 0
$ws.Cells.Item(19, 16).Value = 0
$ws.Cells.Item(19, 17).Value = 17.26722148500747

# Row 20
$ws.Cells.Item(20, 2).Value = 29.46202403993336
$ws.Cells.Item(20, 3).Value = 22.66792729158188
$ws.Cells.Item(20, 4).Value = 8.721716349729293
$ws.Cells.Item(20, 5).Value = 34.45837945924719
$ws.Cells.Item(20, 7).Value = 2.054192750152489
$ws.Cells.Item(20, 8).Value = 4.011724603069078
$ws.Cells.Item(20, 9).Value = 3.842600891627849
$ws.Cells.Item(20, 15).Value = 0
$ws.Cells.Item(20, 16).Value = 0
$ws.Cells.Item(20, 17).Value = 19.72416384284029

# Row 21
$ws.Cells.Item(21, 2).Value = 31.26348911638442
$ws.Cells.Item(21, 3).Value = 23.97662116561174
$ws.Cells.Item(21, 4).Value = 9.286735070701999
$ws.Cells.Item(21, 5).Value = 37.33747205733759
$ws.Cells.Item(21, 7).Value = 2.046419411602481
$ws.Cells.Item(21, 8).Value = 4.398960959570072
$ws.Cells.Item(21, 9).Value = 4.147939749903516
$ws.Cells.Item(21, 15).Value = 0
$ws.Cells.Item(21, 16).Value = 0
$ws.Cells.Item(21, 17).Value = 20.65298020938932

# Row 22
$ws.Cells.Item(22, 2).Value = 32.38962206351611
$ws.Cells.Item(22, 3).Value = 24.75008064270681
$ws.Cells.Item(22, 4).Value = 9.642012757600568
$ws.Cells.Item(22, 5).Value = 38.71064921162453
$ws.Cells.Item(22, 7).Value = 2.041516065056668
$ws.Cells.Item(22, 8).Value = 4.628798883839401
$ws.Cells.Item(22, 9).Value = 4.338982782126293
$ws.Cells.Item(22, 15).Value = 0
$ws.Cells.Item(22, 16).Value = 0
$ws.Cells.Item(22, 17).Value = 21.20491002555745

# Row 23
$ws.Cells.Item(23, 2).Value = 31.79373697220383
$ws.Cells.Item(23, 3).Value = 24.29410116769279
$ws.Cells.Item(23, 4).Value = 9.4537349107962
$ws.Cells.Item(23, 5).Value = 37.98358874188595
$ws.Cells.Item(23, 7).Value = 2.044102609774224
$ws.Cells.Item(23, 8).Value = 4.507043362959229
$ws.Cells.Item(23, 9).Value = 4.235298827237734
$ws.Cells.Item(23, 15).Value = 0
$ws.Cells.Item(23, 16).Value = 0
$ws.Cells.Item(23, 17).Value = 20.9601831962141

# Row 24
$ws.Cells.Item(24, 2).Value = 29.43693827678609
$ws.Cells.Item(24, 3).Value = 22.5803199228251
$ws.Cells.Item(24, 4).Value = 8.713472002555948
$ws.Cells.Item(24, 5).Value = 35.12110194058251
$ws.Cells.Item(24, 7).Value = 2.054046750002974
$ws.Cells.Item(24, 8).Value = 4.038179154857045
$ws.Cells.Item(24, 9).Value = 3.841934490348142
$ws.Cells.Item(24, 15).Value = 0
$ws.Cells.Item(24, 16).Value = 0
$ws.Cells.Item(24, 17).Value = 19.95734931592853

# Row 25
$ws.Cells.Item(25, 2).Value = 26.67983385702391
$ws.Cells.Item(25, 3).Value = 20.60188520522305
$ws.Cells.Item(25, 4).Value = 7.855605411779876
$ws.Cells.Item(25, 5).Value = 31.79903333449962
$ws.Cells.Item(25, 7).Value = 2.065042956790101
$ws.Cells.Item(25, 8).Value = 3.51622204783677
$ws.Cells.Item(25, 9).Value = 3.409616052882697
$ws.Cells.Item(25, 15).Value = 0
$ws.Cells.Item(25, 16).Value = 0
$ws.Cells.Item(25, 17).Value = 18.88658508673164
